$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextCell $ws "D2" "67.054.17"
$ws.Range("E2").Value = "  -1.26%  "

Set-TextCell $ws "D3" "3.579.55"

$ws.Range("E4").Value = "  -0.19%  "

Set-TextCell $ws "D5" "575.09"
$ws.Range("E5").Value = "  -4.60%  "

Set-TextCell $ws "D6" "192.02"
$ws.Range("E6").Value = "  +0.29%  "

Set-TextCell $ws "D7" "3.573.24"
$ws.Range("E7").Value = "  -2.41%  "

Set-TextCell $ws "D8" "0.617"
$ws.Range("E8").Value = "  -1.69%  "

$ws.Range("E9").Value = "  +0.15%  "

Set-TextCell $ws "D10" "0.678"
$ws.Range("E10").Value = "  -5.24%  "

$ws.Range("E11").Value = "  -4.75%  "

Set-TextCell $ws "D12" "56.01"
$ws.Range("E12").Value = "  -5.09%  "

Set-TextCell $ws "D13" "0.0000272"
$ws.Range("E13").Value = "  -3.60%  "

Set-TextCell $ws "D14" "9.84"
$ws.Range("E14").Value = "  -4.44%  "

Set-TextCell $ws "D15" "4.156.28"
$ws.Range("E15").Value = "  -2.58%  "

Set-TextCell $ws "D16" "3.589.96"
$ws.Range("E16").Value = "  -2.48%  "

$ws.Range("E17").Value = "  -1.29%  "

Set-TextCell $ws "D18" "18.37"
$ws.Range("E18").Value = "  -4.27%  "

Set-TextCell $ws "D19" "67.065.18"
$ws.Range("E19").Value = "  -1.16%  "

Set-TextCell $ws "D20" "12.17"
$ws.Range("E20").Value = "  -4.02%  "

Set-TextCell $ws "D21" "1.06"
$ws.Range("E21").Value = "  -6.06%  "

Set-TextCell $ws "D22" "400.56"
$ws.Range("E22").Value = "  -0.78%  "

$ws.Range("E23").Value = "  -7.05%  "

Set-TextCell $ws "D24" "85.77"
$ws.Range("E24").Value = "  -3.64%  "

Set-TextCell $ws "D25" "11.37"
$ws.Range("E25").Value = "  +0.12%  "

Set-TextCell $ws "D26" "2.93"
$ws.Range("E26").Value = "  -3.14%  "

Set-TextCell $ws "D27" "12.46"
$ws.Range("E27").Value = "  -2.92%  "

Set-TextCell $ws "D28" "6.09"
$ws.Range("E28").Value = "  +1.13%  "

Set-TextCell $ws "D29" "3.63"
$ws.Range("E29").Value = "  -2.68%  "

Set-TextCell $ws "D30" "8.96"
$ws.Range("E30").Value = "  -5.44%  "

Set-TextCell $ws "D31" "7.63"
$ws.Range("E31").Value = "  +1.30%  "

Set-TextCell $ws "D32" "31.18"
$ws.Range("E32").Value = "  -3.36%  "

Set-TextCell $ws "D33" "637.97"
$ws.Range("E33").Value = "  +1.91%  "

Set-TextCell $ws "D34" "12.15"
$ws.Range("E34").Value = "  -3.19%  "

Set-TextCell $ws "D35" "0.114"

Set-TextCell $ws "D36" "64.00"
$ws.Range("E36").Value = "  -4.13%  "

Set-TextCell $ws "D37" "42.29"
$ws.Range("E37").Value = "  -9.68%  "

Set-TextCell $ws "D38" "0.400"
$ws.Range("E38").Value = "  -1.47%  "

Set-TextCell $ws "D39" "1.00"
$ws.Range("E39").Value = "  +0.21%  "

Set-TextCell $ws "D40" "0.0₃0774"
$ws.Range("E40").Value = "  -5.13%  "

Set-TextCell $ws "D41" "3.176.68"
$ws.Range("E41").Value = "  +9.68%  "

Set-TextCell $ws "D42" "0.133"
$ws.Range("E42").Value = "  -2.69%  "

$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextCell $ws "D43" "1.00"
$ws.Range("E43").Value = "  -0.14%  "

$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextCell $ws "D44" "2.71"
$ws.Range("E44").Value = "  +4.33%  "

Set-TextCell $ws "D45" "2.98"

Set-TextCell $ws "D46" "0.0416"
$ws.Range("E46").Value = "  -5.22%  "

Set-TextCell $ws "D47" "0.130"
$ws.Range("E47").Value = "  -5.65%  "

$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextCell $ws "D48" "3.08"
$ws.Range("E48").Value = "  +0.70%  "

$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextCell $ws "D49" "2.62"
$ws.Range("E49").Value = "  -1.37%  "

Set-TextCell $ws "D50" "141.62"
$ws.Range("E50").Value = "  -2.82%  "

Set-TextCell $ws "D51" "8.58"
$ws.Range("E51").Value = "  -5.91%  "
